# Optimize-Me.xlsx edit script
# Renames/refreshes the "Parameters" sheet documentation (min_validation_rounds /
# min_models_submitting_per_round notes), widens the Notes column, wraps the new
# long note, resizes row 8, and switches the active tab to Parameters at C23.

$wb = $excel.ActiveWorkbook
$wsModels = $wb.Worksheets.Item("Models")
$wsParams = $wb.Worksheets.Item("Parameters")

# --- Update the "Notes" text for min_validation_rounds (row 7) and
#     min_models_submitting_per_round (row 8). offset_step/roundwindow_step
#     (rows 9/10) keep their original note text.
$wsParams.Range("C7").Value = "Tweak as needed. Higher values gives a larger set of forward OOS validation rounds to evaluate the final models against."
$wsParams.Range("C8").Value = "Tweak as needed. Depending on how spotty model coverage is within a training window, the minimum number of models required affects the actual number of rounds on which is trained. Higher values (up to the total number of models in Optimize-Me) makes for more exclusions of rounds. Lower values makes for more exclusions of models from a training window."
$wsParams.Range("C9").Value = "Tweak as needed. Stepsize for offset sweep. Larger step runs faster; smaller step searches more finely."
$wsParams.Range("C10").Value = "Tweak as needed. Stepsize for windowsize sweep. Larger step runs faster; smaller step searches more finely."

# --- The long note in C8 now wraps and the row grows to fit it.
$wsParams.Range("C8").WrapText = $true
$wsParams.Rows.Item(8).RowHeight = 32

# --- Widen the Notes column so the longer text is readable.
$wsParams.Columns.Item(3).ColumnWidth = 148.33

# --- Models sheet is no longer the active tab; Parameters is, with C23 selected.
$wsParams.Activate()
[void]$wsParams.Range("C23").Select()
